# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.393.57"
$ws.Range("E2").Value = "  +1.54%  "

$ws.Range("D3").Value = "'2.016.04"
$ws.Range("E3").Value = "  +5.14%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'325.23"
$ws.Range("E5").Value = "  +1.74%  "

$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").Value = "'0.5131"
$ws.Range("E7").Value = "  +1.45%  "

$ws.Range("D8").Value = "'0.4262"
$ws.Range("E8").Value = "  +5.68%  "

$ws.Range("D9").Value = "'0.08767"
$ws.Range("E9").Value = "  +5.14%  "

$ws.Range("D10").Value = "'43.48"
$ws.Range("E10").Value = "  +3.28%  "

$ws.Range("D11").Value = "'1.136"
$ws.Range("E11").Value = "  +2.97%  "

$ws.Range("D12").Value = "'24.68"
$ws.Range("E12").Value = "  +3.73%  "

$ws.Range("D13").Value = "'2.011.16"
$ws.Range("E13").Value = "  +5.00%  "

$ws.Range("D14").Value = "'6.608"
$ws.Range("E14").Value = "  +3.29%  "

$ws.Range("D15").Value = "'7.475"
$ws.Range("E15").Value = "  +3.32%  "

$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  -0.39%  "

$ws.Range("D17").Value = "'94.25"
$ws.Range("E17").Value = "  +2.16%  "

$ws.Range("E18").Value = "  +1.61%  "

$ws.Range("D19").Value = "'0.06527"
$ws.Range("E19").Value = "  +0.32%  "

$ws.Range("D20").Value = "'18.91"
$ws.Range("E20").Value = "  +3.82%  "

$ws.Range("E21").Value = "  -0.05%  "

$ws.Range("D22").Value = "'6.218"
$ws.Range("E22").Value = "  +4.58%  "

$ws.Range("D23").Value = "'30.456.89"
$ws.Range("E23").Value = "  +1.68%  "

$ws.Range("D24").Value = "'11.87"
$ws.Range("E24").Value = "  +4.84%  "

$ws.Range("D25").Value = "'2.270"
$ws.Range("E25").Value = "  +3.45%  "

$ws.Range("D26").Value = "'2.248.77"
$ws.Range("E26").Value = "  +4.92%  "

$ws.Range("D27").Value = "'22.50"
$ws.Range("E27").Value = "  +1.89%  "

$ws.Range("D28").Value = "'162.57"
$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("D29").Value = "'2.437"
$ws.Range("E29").Value = "  +5.13%  "

$ws.Range("D30").Value = "'131.31"
$ws.Range("E30").Value = "  +1.40%  "

$ws.Range("D31").Value = "'1.150"
$ws.Range("E31").Value = "  +2.04%  "

$ws.Range("D32").Value = "'0.1054"
$ws.Range("E32").Value = "  +1.49%  "

$ws.Range("D33").Value = "'6.120"
$ws.Range("E33").Value = "  +2.66%  "

$ws.Range("E34").Value = "  +0.36%  "

$ws.Range("E35").Value = "  +15.13%  "

$ws.Range("D36").Value = "'0.02529"
$ws.Range("E36").Value = "  +3.27%  "

$ws.Range("D37").Value = "'5.476"
$ws.Range("E37").Value = "  +1.61%  "

$ws.Range("E38").Value = "  +3.95%  "

$ws.Range("D39").Value = "'12.45"
$ws.Range("E39").Value = "  +9.52%  "

$ws.Range("D40").Value = "'9.159"
$ws.Range("E40").Value = "  +5.37%  "

$ws.Range("D41").Value = "'0.2220"
$ws.Range("E41").Value = "  +2.96%  "

$ws.Range("D42").Value = "'0.6679"
$ws.Range("E42").Value = "  +3.14%  "

$ws.Range("D43").Value = "'1.236"
$ws.Range("E43").Value = "  +1.43%  "

$ws.Range("D44").Value = "'0.9996"
$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("D45").Value = "'13.77"
$ws.Range("E45").Value = "  +2.43%  "

$ws.Range("D46").Value = "'0.6186"
$ws.Range("E46").Value = "  +2.44%  "

$ws.Range("E47").Value = "  -0.74%  "

$ws.Range("D48").Value = "'3.671"
$ws.Range("E48").Value = "  +0.95%  "

$ws.Range("D49").Value = "'1.263"
$ws.Range("E49").Value = "  +4.50%  "

$ws.Range("D50").Value = "'124.91"
$ws.Range("E50").Value = "  +2.24%  "

$ws.Range("D51").Value = "'81.20"
$ws.Range("E51").Value = "  +2.91%  "

